$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.366.74"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").Value = "1.869.68"
$ws.Range("E3").Value = "  +0.08%  "

$c = $ws.Range("D4")
$c.Value = "'1.001"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$c = $ws.Range("D5")
$c.Value = "'329.95"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.37%  "

$ws.Range("E6").Value = "  -0.06%  "

$c = $ws.Range("D7")
$c.Value = "'0.4605"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -2.08%  "

$c = $ws.Range("D8")
$c.Value = "'0.4020"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +2.49%  "

$c = $ws.Range("D9")
$c.Value = "'47.83"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.24%  "

$c = $ws.Range("D10")
$c.Value = "'0.07858"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.72%  "

$c = $ws.Range("D11")
$c.Value = "'0.9856"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.84%  "

$c = $ws.Range("D12")
$c.Value = "'21.30"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -2.03%  "

$ws.Range("D13").Value = "1.868.83"
$ws.Range("E13").Value = "  +0.29%  "

$c = $ws.Range("D14")
$c.Value = "'5.856"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.23%  "

$c = $ws.Range("D15")
$c.Value = "'6.993"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.77%  "

$ws.Range("E16").Value = "  -0.13%  "

$c = $ws.Range("D17")
$c.Value = "'88.17"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -3.18%  "

$c = $ws.Range("D18")
$c.Value = "'0.06547"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.54%  "

$c = $ws.Range("D19")
$c.Value = "'0.00001018"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.35%  "

$c = $ws.Range("D20")
$c.Value = "'17.17"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.51%  "

$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").Value = "28.349.83"
$ws.Range("E22").Value = "  +0.25%  "

$c = $ws.Range("D23")
$c.Value = "'5.337"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.91%  "

$ws.Range("E24").Value = "  -1.73%  "

$ws.Range("E25").Value = "  -1.76%  "

$ws.Range("D26").Value = "2.094.43"
$ws.Range("E26").Value = "  +0.23%  "

$c = $ws.Range("D27")
$c.Value = "'157.70"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.25%  "

$c = $ws.Range("D28")
$c.Value = "'19.34"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -2.26%  "

$c = $ws.Range("D29")
$c.Value = "'2.061"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -3.77%  "

$c = $ws.Range("D30")
$c.Value = "'5.283"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.65%  "

$c = $ws.Range("D31")
$c.Value = "'117.46"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.11%  "

$c = $ws.Range("D32")
$c.Value = "'0.9565"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.91%  "

$c = $ws.Range("D33")
$c.Value = "'0.09325"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.85%  "

$c = $ws.Range("D34")
$c.Value = "'3.587"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.12%  "

$c = $ws.Range("D35")
$c.Value = "'1.386"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.91%  "

$c = $ws.Range("D36")
$c.Value = "'5.234"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.09%  "

$c = $ws.Range("D37")
$c.Value = "'0.06033"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.85%  "

$c = $ws.Range("D38")
$c.Value = "'0.02205"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.03%  "

$c = $ws.Range("D39")
$c.Value = "'8.286"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.59%  "

$c = $ws.Range("D40")
$c.Value = "'1.157"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.76%  "

$c = $ws.Range("D41")
$c.Value = "'1.000"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "

$c = $ws.Range("D42")
$c.Value = "'0.5755"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.47%  "

$ws.Range("E43").Value = "  -3.78%  "

$c = $ws.Range("D44")
$c.Value = "'10.01"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.37%  "

$c = $ws.Range("D45")
$c.Value = "'1.250"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -3.96%  "

$c = $ws.Range("D46")
$c.Value = "'2.297"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +14.90%  "

$ws.Range("E47").Value = "  -3.29%  "

$c = $ws.Range("D48")
$c.Value = "'11.83"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.84%  "

$c = $ws.Range("D49")
$c.Value = "'0.07146"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.60%  "

$c = $ws.Range("D50")
$c.Value = "'1.884"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.04%  "

$c = $ws.Range("D51")
$c.Value = "'110.62"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.16%  "
